$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 218.66667
$ws.Range("I4").Value = 181.75
$ws.Range("J4").Value = 514
$ws.Range("K4").Value = 181.75
$ws.Range("L4").Value = 514
$ws.Range("M4").Value = -67.75
$ws.Range("N4").Value = -742
$ws.Range("H33").Value = 3496608
$ws.Range("I33").Value = 112.210526
$ws.Range("J33").Value = 12987097
$ws.Range("K33").Value = 112.210526
$ws.Range("L33").Value = 12987097
$ws.Range("M33").Value = 116.789474
$ws.Range("N33").Value = -12987555
$ws.Range("H40").Value = 1933.9286
$ws.Range("I40").Value = 1955.4348
$ws.Range("J40").Value = 1835
$ws.Range("K40").Value = 1955.4348
$ws.Range("L40").Value = 1835
$ws.Range("M40").Value = -1780.4348
$ws.Range("N40").Value = -2185
$ws.Range("H106").Value = 83336240
$ws.Range("I106").Value = 37040870
$ws.Range("J106").Value = 142858860
$ws.Range("K106").Value = 37040870
$ws.Range("L106").Value = 142858860
$ws.Range("M106").Value = -37040239
$ws.Range("N106").Value = -142860122
$ws.Range("H116").Value = 6536.364
$ws.Range("I116").Value = 9122.691999999999
$ws.Range("J116").Value = 2800.5557
$ws.Range("K116").Value = 9122.691999999999
$ws.Range("L116").Value = 2800.5557
$ws.Range("M116").Value = -5680.691999999999
$ws.Range("N116").Value = -9684.555700000001
$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178
$ws.Range("H141").Value = 1286.9445
$ws.Range("I141").Value = 1410
$ws.Range("J141").Value = 302.5
$ws.Range("K141").Value = 4230
$ws.Range("L141").Value = 907.5
$ws.Range("M141").Value = 950
$ws.Range("N141").Value = -11267.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5333
$ws.Range("I26").Value = 4799.6
$ws.Range("K26").Value = 4799.6
$ws.Range("M26").Value = -4469.6
$ws.Range("H32").Value = 5226.885
$ws.Range("I32").Value = 3999.7627
$ws.Range("J32").Value = 9037.421
$ws.Range("K32").Value = 3999.7627
$ws.Range("L32").Value = 9037.421
$ws.Range("M32").Value = -3712.7627
$ws.Range("N32").Value = -9611.421
$ws.Range("H74").Value = 1260.7916
$ws.Range("I74").Value = 1145.1875
$ws.Range("J74").Value = 1492
$ws.Range("K74").Value = 1145.1875
$ws.Range("L74").Value = 1492
$ws.Range("M74").Value = -271.1875
$ws.Range("N74").Value = -3240
$ws.Range("H77").Value = 1260.7916
$ws.Range("I77").Value = 1145.1875
$ws.Range("J77").Value = 1492
$ws.Range("K77").Value = 5725.9375
$ws.Range("L77").Value = 7460
$ws.Range("M77").Value = -1357.9375
$ws.Range("N77").Value = -16196
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 335466.66
$ws.Range("I24").Value = 335466.66
$ws.Range("K24").Value = 335466.66
$ws.Range("M24").Value = -335231.66
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4049861.2
$ws.Range("I16").Value = 12821295
$ws.Range("J16").Value = 1507
$ws.Range("K16").Value = 12821295
$ws.Range("L16").Value = 1507
$ws.Range("M16").Value = -12821008
$ws.Range("N16").Value = -2081
$ws.Range("H31").Value = 5046.1753
$ws.Range("I31").Value = 1688.5555
$ws.Range("J31").Value = 10802.096
$ws.Range("K31").Value = 1688.5555
$ws.Range("L31").Value = 10802.096
$ws.Range("M31").Value = -1393.5555
$ws.Range("N31").Value = -11392.096
$ws.Range("H34").Value = 5046.1753
$ws.Range("I34").Value = 1688.5555
$ws.Range("J34").Value = 10802.096
$ws.Range("K34").Value = 1688.5555
$ws.Range("L34").Value = 10802.096
$ws.Range("M34").Value = -1486.5555
$ws.Range("N34").Value = -11206.096
$ws.Range("H113").Value = 4049861.2
$ws.Range("I113").Value = 12821295
$ws.Range("J113").Value = 1507
$ws.Range("K113").Value = 12821295
$ws.Range("L113").Value = 1507
$ws.Range("M113").Value = -12819125
$ws.Range("N113").Value = -5847
$ws.Range("H132").Value = 1483.9265
$ws.Range("I132").Value = 1017.1404
$ws.Range("J132").Value = 3902.7273
$ws.Range("K132").Value = 3051.4212
$ws.Range("L132").Value = 11708.1819
$ws.Range("M132").Value = -521.4211999999998
$ws.Range("N132").Value = -16768.1819
$ws.Range("H134").Value = 3147.862
$ws.Range("I134").Value = 4200.353
$ws.Range("J134").Value = 1656.8334
$ws.Range("K134").Value = 12601.059
$ws.Range("L134").Value = 4970.5002
$ws.Range("M134").Value = -10066.059
$ws.Range("N134").Value = -10040.5002
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2740
$ws.Range("I25").Value = 1200
$ws.Range("J25").Value = 3125
$ws.Range("K25").Value = 3600
$ws.Range("L25").Value = 9375
$ws.Range("N25").Value = -9713
$ws.Range("M25").Value = -3431
$ws.Range("H30").Value = 2740
$ws.Range("I30").Value = 1200
$ws.Range("J30").Value = 3125
$ws.Range("K30").Value = 3600
$ws.Range("L30").Value = 9375
$ws.Range("N30").Value = -9579
$ws.Range("M30").Value = -3498
$ws.Range("H113").Value = 455123.38
$ws.Range("J113").Value = 1250487.5
$ws.Range("L113").Value = 3751462.5
$ws.Range("N113").Value = -3755802.5
$ws.Range("H131").Value = 2381756.8
$ws.Range("I131").Value = 6250317.5
$ws.Range("J131").Value = 1103.9231
$ws.Range("K131").Value = 18750952.5
$ws.Range("L131").Value = 3311.7693
$ws.Range("M131").Value = -18745912.5
$ws.Range("N131").Value = -13391.7693
$ws.Range("H136").Value = 24010.6
$ws.Range("J136").Value = 5005.75
$ws.Range("L136").Value = 15017.25
$ws.Range("N136").Value = -25217.25
$ws.Range("H137").Value = 23820392
$ws.Range("I137").Value = 16865.715
$ws.Range("J137").Value = 47623916
$ws.Range("K137").Value = 50597.145
$ws.Range("L137").Value = 142871748
$ws.Range("M137").Value = -45497.145
$ws.Range("N137").Value = -142881948
$ws.Range("H138").Value = 9986.412
$ws.Range("I138").Value = 14039.818
$ws.Range("J138").Value = 2555.1667
$ws.Range("K138").Value = 42119.454
$ws.Range("L138").Value = 7665.500100000001
$ws.Range("M138").Value = -36979.454
$ws.Range("N138").Value = -17945.5001
$ws.Range("H139").Value = 6478.077
$ws.Range("I139").Value = 21698
$ws.Range("J139").Value = 2854.2856
$ws.Range("K139").Value = 65094
$ws.Range("L139").Value = 8562.856800000001
$ws.Range("M139").Value = -59954
$ws.Range("N139").Value = -18842.8568
$ws.Range("H140").Value = 5337.067
$ws.Range("I140").Value = 5337.067
$ws.Range("K140").Value = 16011.201
$ws.Range("M140").Value = -10831.201
$ws.Range("H141").Value = 140516.38
$ws.Range("I141").Value = 275682.75
$ws.Range("J141").Value = 5350
$ws.Range("K141").Value = 827048.25
$ws.Range("L141").Value = 16050
$ws.Range("M141").Value = -821868.25
$ws.Range("N141").Value = -26410
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 14000000
$ws.Range("I24").Value = 14000000
$ws.Range("K24").Value = 14000000
$ws.Range("M24").Value = -13999827
$ws.Range("H126").Value = 4840.579
$ws.Range("I126").Value = 7665.4707
$ws.Range("J126").Value = 2553.762
$ws.Range("K126").Value = 22996.4121
$ws.Range("L126").Value = 7661.286
$ws.Range("M126").Value = -20526.4121
$ws.Range("N126").Value = -12601.286
$ws.Range("H132").Value = 2937.6191
$ws.Range("I132").Value = 2565.1428
$ws.Range("J132").Value = 3682.5715
$ws.Range("K132").Value = 7695.428400000001
$ws.Range("L132").Value = 11047.7145
$ws.Range("M132").Value = -5165.428400000001
$ws.Range("N132").Value = -16107.7145
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 3274.75
$ws.Range("I107").Value = 3274.75
$ws.Range("K107").Value = 3274.75
$ws.Range("M107").Value = -1354.75
$ws.Range("H132").Value = 8991465
$ws.Range("I132").Value = 11183542
$ws.Range("J132").Value = 3950.8
$ws.Range("K132").Value = 33550626
$ws.Range("L132").Value = 11852.4
$ws.Range("M132").Value = -33548096
$ws.Range("N132").Value = -16912.4
$ws.Range("H136").Value = 6785.174
$ws.Range("I136").Value = 4830.147
$ws.Range("K136").Value = 14490.441
$ws.Range("M136").Value = -11940.441
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H18").Value = 7210
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 7210
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 7210
$ws.Range("M18").ClearContents()
$ws.Range("N18").Value = -7556
$ws.Range("H20").Value = 5003605
$ws.Range("I20").Value = 10000000
$ws.Range("K20").Value = 10000000
$ws.Range("M20").Value = -9999760
$ws.Range("H21").Value = 38588.75
$ws.Range("I21").Value = 7101
$ws.Range("J21").Value = 49084.668
$ws.Range("K21").Value = 7101
$ws.Range("L21").Value = 49084.668
$ws.Range("N21").Value = -49554.668
$ws.Range("M21").Value = -6866
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H35").Value = 38588.75
$ws.Range("I35").Value = 7101
$ws.Range("J35").Value = 49084.668
$ws.Range("K35").Value = 7101
$ws.Range("L35").Value = 49084.668
$ws.Range("N35").Value = -49664.668
$ws.Range("M35").Value = -6811
$ws.Range("H132").Value = 1077.3948
$ws.Range("I132").Value = 731.0909
$ws.Range("J132").Value = 1984.381
$ws.Range("K132").Value = 2193.2727
$ws.Range("L132").Value = 5953.143
$ws.Range("M132").Value = 336.7273
$ws.Range("N132").Value = -11013.143
$ws.Range("H136").Value = 1791.9636
$ws.Range("I136").Value = 1720.8788
$ws.Range("J136").Value = 1898.591
$ws.Range("K136").Value = 5162.636399999999
$ws.Range("L136").Value = 5695.772999999999
$ws.Range("M136").Value = -2612.636399999999
$ws.Range("N136").Value = -10795.773
